$d = $word.ActiveDocument

# The document's visible content lives entirely in the first paragraph;
# it is followed by a handful of empty trailing paragraphs. We need to:
#   1. Append "/ limCrypt / cinder-block-and-plywood shantytowns" right
#      after the existing text (before the _GoBack bookmark) in that
#      first paragraph, using the same run formatting already in use
#      throughout the document (Untitled Serif, 20pt/sz 40, ligatures).
#   2. Remove the trailing empty paragraphs entirely.

$firstPara = $d.Paragraphs.Item(1)

# Build a zero-length range that sits right before the paragraph mark
# (and therefore before the bookmark, which is the last thing in the
# paragraph) so inserted text lands inside the paragraph, ahead of the
# bookmark markers.
$insertionPoint = $firstPara.Range
$insertionPoint.MoveEnd(1, -1) | Out-Null
$insertionPoint.Collapse(0)

# Match the run formatting used throughout the rest of the document.
$insertionPoint.Font.Name = "Untitled Serif"
$insertionPoint.Font.Size = 20

$insertionPoint.InsertBefore("/ limCrypt / cinder-block-and-plywood shantytowns")

# Drop the now-stale empty paragraphs that trailed the content
# paragraph (everything from the end of paragraph 1 to the end of the
# document).
$firstPara = $d.Paragraphs.Item(1)
$tail = $d.Range($firstPara.Range.End, $d.Content.End)
if ($tail.Start -lt $tail.End) {
    $tail.Delete()
}

Write-Host "Paragraph count:" $d.Paragraphs.Count
Write-Host "Final text:" $d.Paragraphs.Item(1).Range.Text
